$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC23_Verify_UserRegistration")
$ws2 = $wb.Worksheets.Item("Testdata")

# Insert a new row at position 27 (TC23 sheet), shifting rows 27-29 down to 28-30.
$ws1.Rows.Item(27).Insert()

# New row 27 becomes a WAIT step (matching the other WAIT rows in the sheet).
$ws1.Range("B27").Value = "WAIT"

# Match the thin-border look used by every other data row on this sheet.
$ws1.Range("A27:E27").Borders.LineStyle = 1

# Restore cursor/selection state: sheet2 selection moves to L17, sheet1 selection to B27.
[void]$ws2.Range("L17").Select()
[void]$ws1.Activate()
[void]$ws1.Range("B27").Select()
